$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Kenneth A.%Egol%NULL%2,         Sanjit R.%Konda%NULL%2,         Mackenzie L.%Bird%NULL%2,         Nicket%Dedhia%NULL%2,         Emma K.%Landes%NULL%2,         Rachel A.%Ranson%NULL%2,         Sara J.%Solasz%NULL%2,         Vinay K.%Aggarwal%NULL%2,         Joseph A.%Bosco%NULL%2,         David L.%Furgiuele%NULL%2,         Abhishek%Ganta%NULL%2,         Jason%Gould%NULL%2,         Thomas R.%Lyon%NULL%2,         Toni M.%McLaurin%NULL%2,         Nirmal C.%Tejwani%NULL%2,         Joseph D.%Zuckerman%NULL%2,         Philipp%Leucht%NULL%2]'
$ws.Range("E3").Value = '[Drake G.%LeBrun%NULL%2,         Maxwell A.%Konnaris%NULL%2,         Gregory C.%Ghahramani%NULL%2,         Ajay%Premkumar%NULL%2,         Chris J.%DeFrancesco%NULL%2,         Jordan A.%Gruskay%NULL%2,         Aleksey%Dvorzhinskiy%NULL%2,         Milan S.%Sandhu%NULL%2,         Elan M.%Goldwyn%NULL%2,         Christopher L.%Mendias%NULL%2,         William M.%Ricci%NULL%2]'
$ws.Range("E4").Value = '[Josep Maria%Muñoz Vives%NULL%1,         Montsant%Jornet-Gibert%NULL%2,         Montsant%Jornet-Gibert%NULL%0,         J.%Cámara-Cabrera%NULL%2,         J.%Cámara-Cabrera%NULL%0,         Pedro L.%Esteban%NULL%2,         Pedro L.%Esteban%NULL%0,         Laia%Brunet%NULL%2,         Laia%Brunet%NULL%0,         Luis%Delgado-Flores%NULL%2,         Luis%Delgado-Flores%NULL%0,         P.%Camacho-Carrasco%NULL%2,         P.%Camacho-Carrasco%NULL%0,         P.%Torner%NULL%2,         P.%Torner%NULL%0,         Francesc%Marcano-Fernández%NULL%2,         Francesc%Marcano-Fernández%NULL%0,         NULL%NULL%NULL%0,         NULL%NULL%NULL%0]'
$ws.Range("E5").Value = '[Zoe B.%Cheung%zoe.cheung@mountsinai.org%1,         David A.%Forsh%NULL%1]'
$ws.Range("C6").Value = '"IMPACT-Scot report on COVID-19 and hip fractures"'
$ws.Range("D6").Value = '"Aims The primary aim was to assess the independent influence of coronavirus disease (COVID-19) on 30-day mortality for patients with a hip fracture.
 The secondary aims were to determine whether: 1) there were clinical predictors of COVID-19 status; and 2) whether social lockdown influenced the incidence and epidemiology of hip fractures.
 Methods A national multicentre retrospective study was conducted of all patients presenting to six trauma centres or units with a hip fracture over a 46-day period (23 days pre- and 23 days post-lockdown).
 Patient demographics, type of residence, place of injury, presentation blood tests, Nottingham Hip Fracture Score, time to surgery, operation, American Society of Anesthesiologists (ASA) grade, anaesthetic, length of stay, COVID-19 status, and 30-day mortality were recorded.
 Results Of 317 patients with acute hip fracture, 27 (8.5%) had a positive COVID-19 test.
 Only seven (26%) had suggestive symptoms on admission.
 COVID-19-positive patients had a significantly lower 30-day survival compared to those without COVID-19 (64.5%, 95% confidence interval (CI) 45.7 to 83.3 vs 91.7%, 95% CI 88.2 to 94.8; p &lt; 0.001).
 COVID-19 was independently associated with increased 30-day mortality risk adjusting for: 1) age, sex, type of residence (hazard ratio (HR) 2.93; p = 0.008); 2) Nottingham Hip Fracture Score (HR 3.52; p = 0.001); and 3) ASA (HR 3.45; p = 0.004).
 Presentation platelet count predicted subsequent COVID-19 status; a value of &lt; 217 × 109/l was associated with 68% area under the curve (95% CI 58 to 77; p = 0.002) and a sensitivity and specificity of 63%.
 A similar number of patients presented with hip fracture in the 23 days pre-lockdown (n = 160) and 23 days post-lockdown (n = 157) with no significant (all p ≥ 0.130) difference in patient demographics, residence, place of injury, Nottingham Hip Fracture Score, time to surgery, ASA, or management.
 Conclusion COVID-19 was independently associated with an increased 30-day mortality rate for patients with a hip fracture.
 Notably, most patients with hip fracture and COVID-19 lacked suggestive symptoms at presentation.
 Platelet count was an indicator of risk of COVID-19 infection.
 These findings have implications for the management of hip fractures, in particular the need for COVID-19 testing.
 Cite this article: Bone Joint J 2020;102-B(9):1219–1228. '
$ws.Range("E6").Value = '[Andrew J.%Hall%xref no email%2,  Nicholas D.%Clement%xref no email%2,  Luke%Farrow%xref no email%2,  Alasdair M. J.%MacLullich%xref no email%2,  Graham F.%Dall%xref no email%2,  Chloe E. H.%Scott%xref no email%2,  Paul J.%Jenkins%xref no email%2,  Timothy O.%White%xref no email%2,  Andrew D.%Duckworth%xref no email%2]'
$ws.Range("F6").Value = '10.1302/0301-620X.102B9.BJJ-2020-1100.R1'
$ws.Range("G6").Value = 'CROSSREF'
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = '2023-05-26'
$ws.Range("H6").ClearFormats()
$ws.Range("C7").Value = '"The effects of COVID-19 on perioperative morbidity and mortality in patients with hip fractures"'
$ws.Range("D7").Value = '"Aims During the COVID-19 pandemic, many patients continue to require urgent surgery for hip fractures.
 However, the impact of COVID-19 on perioperative outcomes in these high-risk patients remains unknown.
 The objectives of this study were to establish the effects of COVID-19 on perioperative morbidity and mortality, and determine any risk factors for increased mortality in patients with COVID-19 undergoing hip fracture surgery.
 Methods This multicentre cohort study included 340 COVID-19-negative patients versus 82 COVID-19-positive patients undergoing surgical treatment for hip fractures across nine NHS hospitals in Greater London, UK.
 Patients in both treatment groups were comparable for age, sex, body mass index, fracture configuration, and type of surgery performed.
 Predefined perioperative outcomes were recorded within a 30-day postoperative period.
 Univariate and multivariate analysis were used to identify risk factors associated with increased risk of mortality.
 Results COVID-19-positive patients had increased postoperative mortality rates (30.5% (25/82) vs 10.3% (35/340) respectively, p &lt; 0.001) compared to COVID-19-negative patients.
 Risk factors for increased mortality in patients with COVID-19 undergoing surgery included positive smoking status (hazard ratio (HR) 15.4 (95% confidence interval (CI) 4.55 to 52.2; p &lt; 0.001) and greater than three comorbidities (HR 13.5 (95% CI 2.82 to 66.0, p &lt; 0.001).
 COVID-19-positive patients had increased risk of postoperative complications (89.0% (73/82) vs 35.0% (119/340) respectively; p &lt; 0.001), more critical care unit admissions (61.0% (50/82) vs 18.2% (62/340) respectively; p &lt; 0.001), and increased length of hospital stay (mean 13.8 days (SD 4.6) vs 6.7 days (SD 2.5) respectively; p &lt; 0.001), compared to COVID-19-negative patients.
 Conclusion Hip fracture surgery in COVID-19-positive patients was associated with increased length of hospital stay, more admissions to the critical care unit, higher risk of perioperative complications, and increased mortality rates compared to COVID-19-negative patients.
 Risk factors for increased mortality in patients with COVID-19 undergoing surgery included positive smoking status and multiple (greater than three) comorbidities.
 Cite this article: Bone Joint J 2020;102-B(9):1136–1145. '
$ws.Range("E7").Value = '[Babar%Kayani%xref no email%2,  Elliot%Onochie%xref no email%2,  Vijay%Patil%xref no email%2,  Fahima%Begum%xref no email%2,  Rory%Cuthbert%xref no email%2,  David%Ferguson%xref no email%2,  Jagmeet S.%Bhamra%xref no email%2,  Aadhar%Sharma%xref no email%2,  Peter%Bates%xref no email%2,  Fares S.%Haddad%xref no email%2]'
$ws.Range("F7").Value = '10.1302/0301-620X.102B9.BJJ-2020-1127.R1'
$ws.Range("G7").Value = 'CROSSREF'
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = '2023-05-26'
$ws.Range("H7").ClearFormats()
